$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("Q2").Value = 2.06
$ws.Range("R2").Value = 1.84

# Row 4 updates
$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 2.38
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 9.5
$ws.Range("W4").Value = 9.5
$ws.Range("AC4").Value = 9.5
$ws.Range("AE4").Value = 13
$ws.Range("AL4").Value = 19

# Row 6 updates
$ws.Range("G6").Value = 2.1
$ws.Range("I6").Value = 2.72
$ws.Range("J6").Value = 2.47
$ws.Range("K6").Value = 2.7
$ws.Range("L6").Value = 3
$ws.Range("P6").Value = 6.5
$ws.Range("Q6").Value = 1.27
$ws.Range("R6").Value = 3.4
$ws.Range("S6").Value = 1.17
$ws.Range("T6").Value = 4.45
$ws.Range("U6").Value = 1.29
$ws.Range("V6").Value = 3.3
$ws.Range("W6").Value = 19
$ws.Range("X6").Value = 18.5
$ws.Range("Y6").Value = 10.25
$ws.Range("Z6").Value = 26
$ws.Range("AD6").Value = 10.75
$ws.Range("AJ6").Value = 12
$ws.Range("AK6").Value = 40
$ws.Range("AL6").Value = 19
$ws.Range("AN6").Value = 5.2
$ws.Range("AO6").Value = 10.25
$ws.Range("AP6").Value = 12
$ws.Range("AQ6").Value = 29
$ws.Range("AS6").Value = 75
$ws.Range("AT6").Value = 4.45
$ws.Range("AV6").Value = 26
$ws.Range("AX6").Value = 13.5
$ws.Range("AY6").Value = 13.5
$ws.Range("AZ6").Value = 45
$ws.Range("BA6").Value = 45
$ws.Range("BC6").Value = 250

$wb.Save()
